# Updated cryptos list on Tue Sep 10 23:43:15 UTC 2024 with GitHub Actions
#
# Writes each changed cell as plain TEXT (matching the source data's
# inline-string cells, e.g. "518.53", "10.60", "  +0.35%  ") rather than
# letting Excel auto-coerce numeric-looking strings into numbers, and
# restores the cell's style afterwards so no stray number-format/style
# gets left behind on the cell.
function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Price (D) / Volume(1h) (E) updates -------------------------
Set-TextCell $ws "D2"  "57.402.21"
Set-TextCell $ws "E2"  "  +0.38%  "

Set-TextCell $ws "D3"  "2.366.44"
Set-TextCell $ws "E3"  "  -0.02%  "

Set-TextCell $ws "D5"  "518.45"
Set-TextCell $ws "E5"  "  -0.31%  "

Set-TextCell $ws "D6"  "135.66"
Set-TextCell $ws "E6"  "  +0.25%  "

Set-TextCell $ws "E7"  "  -0.41%  "

Set-TextCell $ws "E9"  "  -1.01%  "

Set-TextCell $ws "D10" "5.53"
Set-TextCell $ws "E10" "  +5.69%  "

Set-TextCell $ws "E11" "  -0.95%  "

Set-TextCell $ws "E12" "  -0.08%  "

Set-TextCell $ws "D13" "24.35"
Set-TextCell $ws "E13" "  +1.81%  "

Set-TextCell $ws "D14" "2.792.88"
Set-TextCell $ws "E14" "  +0.22%  "

Set-TextCell $ws "D15" "57.402.99"
Set-TextCell $ws "E15" "  +0.58%  "

Set-TextCell $ws "E16" "  +0.16%  "

Set-TextCell $ws "D17" "2.377.08"
Set-TextCell $ws "E17" "  +0.58%  "

Set-TextCell $ws "D18" "10.59"
Set-TextCell $ws "E18" "  +0.21%  "

Set-TextCell $ws "D19" "330.13"
Set-TextCell $ws "E19" "  +2.10%  "

Set-TextCell $ws "E20" "  -0.75%  "

Set-TextCell $ws "E21" "  -0.41%  "

Set-TextCell $ws "E22" "  -0.34%  "

Set-TextCell $ws "D23" "61.47"
Set-TextCell $ws "E23" "  -0.03%  "

Set-TextCell $ws "D24" "8.95"
Set-TextCell $ws "E24" "  +14.76%  "

Set-TextCell $ws "E25" "  +3.52%  "

Set-TextCell $ws "D26" "0.995"
Set-TextCell $ws "E26" "  -0.15%  "

Set-TextCell $ws "E27" "  +10.67%  "

Set-TextCell $ws "E28" "  +0.63%  "

# --- Rows 29/30 swap: PancakeSwap <-> Monero ---------------------------
Set-TextCell $ws "B29" "Monero"
Set-TextCell $ws "C29" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D29" "167.07"
Set-TextCell $ws "E29" "  -2.39%  "

Set-TextCell $ws "B30" "PancakeSwap"
Set-TextCell $ws "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D30" "1.70"
Set-TextCell $ws "E30" "  +0.99%  "

Set-TextCell $ws "D31" "6.27"
Set-TextCell $ws "E31" "  -0.41%  "

Set-TextCell $ws "D32" "18.61"
Set-TextCell $ws "E32" "  +1.05%  "

Set-TextCell $ws "D34" "1.30"
Set-TextCell $ws "E34" "  +3.24%  "

Set-TextCell $ws "E35" "  -0.51%  "

Set-TextCell $ws "D36" "0.919"
Set-TextCell $ws "E36" "  -3.78%  "

Set-TextCell $ws "E37" "  +0.52%  "

Set-TextCell $ws "E38" "  +6.15%  "

Set-TextCell $ws "D39" "38.90"
Set-TextCell $ws "E39" "  +3.57%  "

Set-TextCell $ws "D40" "150.48"
Set-TextCell $ws "E40" "  +6.92%  "

Set-TextCell $ws "E41" "  +1.13%  "

Set-TextCell $ws "D43" "289.37"
Set-TextCell $ws "E43" "  +3.79%  "

Set-TextCell $ws "D44" "5.29"
Set-TextCell $ws "E44" "  +2.22%  "

Set-TextCell $ws "E45" "  +1.02%  "

Set-TextCell $ws "E46" "  -0.47%  "

Set-TextCell $ws "E47" "  +0.88%  "

Set-TextCell $ws "D48" "0.389"
Set-TextCell $ws "E48" "  +1.56%  "

# --- Rows 49/50 swap: InjectiveProtocol <-> EnergySwap -----------------
Set-TextCell $ws "B49" "EnergySwap"
Set-TextCell $ws "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D49" "17.78"
Set-TextCell $ws "E49" "  +4.48%  "

Set-TextCell $ws "B50" "InjectiveProtocol"
Set-TextCell $ws "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D50" "18.19"
Set-TextCell $ws "E50" "  +5.03%  "

Write-Host "Edit applied"
